$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4 (shifts existing rows 4-18 down to 5-19)
$ws.Rows.Item(4).Insert()

# Fill in the new SeniorCitizen row
$ws.Cells.Item(4, 1).Value = "SeniorCitizen"
$ws.Cells.Item(4, 2).Value = 7043
$ws.Cells.Item(4, 3).Value = 2
$ws.Cells.Item(4, 4).Value = "No"
$ws.Cells.Item(4, 5).Value = 5901
$ws.Cells.Item(4, 6).Value = 83.78531875621185

# Copy the style of column A from row 3 (gender) to the new row 4 cell
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(4, 1).PasteSpecial(-4122) | Out-Null
